$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2841173
$ws.Range("J17").Value = 2841173
$ws.Range("L17").Value = 8523519
$ws.Range("N17").Value = -8523855
$ws.Range("H40").Value = 2767
$ws.Range("J40").Value = 4000
$ws.Range("L40").Value = 4000
$ws.Range("N40").Value = -4350
$ws.Range("H70").Value = 7411474.5
$ws.Range("I70").Value = 25006780
$ws.Range("J70").Value = 2925.0527
$ws.Range("K70").Value = 75020340
$ws.Range("L70").Value = 8775.158100000001
$ws.Range("M70").Value = -75020070
$ws.Range("N70").Value = -9315.158100000001
$ws.Range("H73").Value = 7411474.5
$ws.Range("I73").Value = 25006780
$ws.Range("J73").Value = 2925.0527
$ws.Range("K73").Value = 75020340
$ws.Range("L73").Value = 8775.158100000001
$ws.Range("M73").Value = -75019404
$ws.Range("N73").Value = -10647.1581
$ws.Range("H103").Value = 3309.4546
$ws.Range("I103").Value = 1219.3636
$ws.Range("J103").Value = 5399.5454
$ws.Range("K103").Value = 3658.0908
$ws.Range("L103").Value = 16198.6362
$ws.Range("M103").Value = -3072.0908
$ws.Range("N103").Value = -17370.6362
$ws.Range("H112").Value = 49115.52
$ws.Range("I112").Value = 1599
$ws.Range("J112").Value = 55595.047
$ws.Range("K112").Value = 4797
$ws.Range("L112").Value = 166785.141
$ws.Range("M112").Value = -3689
$ws.Range("N112").Value = -169001.141
$ws.Range("H121").Value = 973.8182
$ws.Range("J121").Value = 973.8182
$ws.Range("L121").Value = 2921.4546
$ws.Range("N121").Value = -6415.4546
$ws.Range("H125").Value = 4047.2222
$ws.Range("J125").Value = 4474.5
$ws.Range("L125").Value = 40270.5
$ws.Range("N125").Value = -45190.5
$ws.Range("H131").Value = 4679.875
$ws.Range("I131").Value = 2487.8
$ws.Range("K131").Value = 7463.400000000001
$ws.Range("M131").Value = -2423.400000000001
$ws.Range("H135").Value = 3276.111
$ws.Range("I135").Value = 3603.7144
$ws.Range("K135").Value = 32433.4296
$ws.Range("M135").Value = -29898.4296
$ws.Range("H138").Value = 8444.944
$ws.Range("J138").Value = 7951.6
$ws.Range("L138").Value = 23854.8
$ws.Range("N138").Value = -34134.8
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3076.72
$ws.Range("I61").Value = 2707.4666
$ws.Range("K61").Value = 2707.4666
$ws.Range("M61").Value = -2495.4666
$ws.Range("H74").Value = 26392680
$ws.Range("I74").Value = 27858662
$ws.Range("K74").Value = 27858662
$ws.Range("M74").Value = -27857788
$ws.Range("H77").Value = 26392680
$ws.Range("I77").Value = 27858662
$ws.Range("K77").Value = 139293310
$ws.Range("M77").Value = -139288942
$ws.Range("H122").Value = 4127.7334
$ws.Range("I122").Value = 3141.4211
$ws.Range("J122").Value = 5831.364
$ws.Range("K122").Value = 9424.263300000001
$ws.Range("L122").Value = 17494.092
$ws.Range("M122").Value = -6974.263300000001
$ws.Range("N122").Value = -22394.092
$ws.Range("H136").Value = 3076.72
$ws.Range("I136").Value = 2707.4666
$ws.Range("K136").Value = 8122.399800000001
$ws.Range("M136").Value = -5572.399800000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 7921.4
$ws.Range("I94").Value = 8945.857
$ws.Range("J94").Value = 4335.8
$ws.Range("K94").Value = 8945.857
$ws.Range("L94").Value = 4335.8
$ws.Range("M94").Value = -8494.857
$ws.Range("N94").Value = -5237.8
$ws.Range("H107").Value = 2496.1667
$ws.Range("I107").Value = 2795.8
$ws.Range("J107").Value = 998
$ws.Range("K107").Value = 2795.8
$ws.Range("L107").Value = 998
$ws.Range("M107").Value = -875.8000000000002
$ws.Range("N107").Value = -4838
$ws.Range("H134").Value = 2750.56
$ws.Range("I134").Value = 2671.5
$ws.Range("J134").Value = 3330.3333
$ws.Range("K134").Value = 8014.5
$ws.Range("L134").Value = 9990.999899999999
$ws.Range("M134").Value = -5479.5
$ws.Range("N134").Value = -15060.9999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2004.2297
$ws.Range("I31").Value = 862.25
$ws.Range("J31").Value = 2552.38
$ws.Range("K31").Value = 862.25
$ws.Range("L31").Value = 2552.38
$ws.Range("M31").Value = -567.25
$ws.Range("N31").Value = -3142.38
$ws.Range("H34").Value = 2004.2297
$ws.Range("I34").Value = 862.25
$ws.Range("J34").Value = 2552.38
$ws.Range("K34").Value = 862.25
$ws.Range("L34").Value = 2552.38
$ws.Range("M34").Value = -660.25
$ws.Range("N34").Value = -2956.38
$ws.Range("H58").Value = 5734.7856
$ws.Range("I58").Value = 7101.087
$ws.Range("K58").Value = 7101.087
$ws.Range("M58").Value = -6898.087
$ws.Range("H134").Value = 1331572.8
$ws.Range("I134").Value = 1895230.8
$ws.Range("K134").Value = 5685692.4
$ws.Range("M134").Value = -5683157.4
$ws.Range("H136").Value = 5734.7856
$ws.Range("I136").Value = 7101.087
$ws.Range("K136").Value = 21303.261
$ws.Range("M136").Value = -18753.261
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 200.66667
$ws.Range("J2").Value = 284.33334
$ws.Range("L2").Value = 1706.00004
$ws.Range("N2").Value = -1932.00004
$ws.Range("H34").Value = 1736266.2
$ws.Range("I34").Value = 2780149.2
$ws.Range("J34").Value = 170441.67
$ws.Range("K34").Value = 8340447.600000001
$ws.Range("L34").Value = 511325.01
$ws.Range("M34").Value = -8340363.600000001
$ws.Range("N34").Value = -511493.01
$ws.Range("H38").Value = 1443.8096
$ws.Range("J38").Value = 2033.2858
$ws.Range("L38").Value = 6099.857400000001
$ws.Range("N38").Value = -6793.857400000001
$ws.Range("H39").Value = 8380.833000000001
$ws.Range("J39").Value = 9997
$ws.Range("L39").Value = 29991
$ws.Range("N39").Value = -30579
$ws.Range("H68").Value = 7195.1304
$ws.Range("J68").Value = 9919.532999999999
$ws.Range("L68").Value = 29758.599
$ws.Range("N68").Value = -31380.599
$ws.Range("H71").Value = 7195.1304
$ws.Range("J71").Value = 9919.532999999999
$ws.Range("L71").Value = 89275.79699999999
$ws.Range("N71").Value = -97387.79699999999
$ws.Range("H121").Value = 667388.5
$ws.Range("J121").Value = 2138.6667
$ws.Range("L121").Value = 6416.000100000001
$ws.Range("N121").Value = -9036.000100000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5164.147
$ws.Range("I70").Value = 5161.276
$ws.Range("K70").Value = 5161.276
$ws.Range("M70").Value = -4891.276
$ws.Range("H73").Value = 5164.147
$ws.Range("I73").Value = 5161.276
$ws.Range("K73").Value = 5161.276
$ws.Range("M73").Value = -4225.276
$ws.Range("H102").Value = 3894.125
$ws.Range("I102").Value = 3894.125
$ws.Range("K102").Value = 3894.125
$ws.Range("M102").Value = -2272.125
$ws.Range("H126").Value = 8543.666999999999
$ws.Range("I126").Value = 6783.25
$ws.Range("J126").Value = 10304.083
$ws.Range("K126").Value = 20349.75
$ws.Range("L126").Value = 30912.249
$ws.Range("M126").Value = -17879.75
$ws.Range("N126").Value = -35852.249
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5799.2954
$ws.Range("I22").Value = 7020.56
$ws.Range("K22").Value = 7020.56
$ws.Range("M22").Value = -6725.56
$ws.Range("H27").Value = 5799.2954
$ws.Range("I27").Value = 7020.56
$ws.Range("K27").Value = 7020.56
$ws.Range("M27").Value = -6913.56
$ws.Range("H61").Value = 2309.3333
$ws.Range("I61").Value = 1783
$ws.Range("K61").Value = 1783
$ws.Range("M61").Value = -1581
$ws.Range("H113").Value = 2309.3333
$ws.Range("I113").Value = 1783
$ws.Range("K113").Value = 1783
$ws.Range("M113").Value = 387
$ws.Range("H122").Value = 4614.8667
$ws.Range("I122").Value = 4167.091
$ws.Range("K122").Value = 12501.273
$ws.Range("M122").Value = -10051.273
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 20491.25
$ws.Range("I62").Value = 33407.152
$ws.Range("J62").Value = 5227
$ws.Range("K62").Value = 33407.152
$ws.Range("L62").Value = 5227
$ws.Range("M62").Value = -32783.152
$ws.Range("N62").Value = -6475
$ws.Range("H65").Value = 20491.25
$ws.Range("I65").Value = 33407.152
$ws.Range("J65").Value = 5227
$ws.Range("K65").Value = 167035.76
$ws.Range("L65").Value = 26135
$ws.Range("M65").Value = -163915.76
$ws.Range("N65").Value = -32375
$ws.Range("H74").Value = 20463.8
$ws.Range("J74").Value = 20463.8
$ws.Range("L74").Value = 20463.8
$ws.Range("N74").Value = -22335.8
$ws.Range("H77").Value = 20463.8
$ws.Range("J77").Value = 20463.8
$ws.Range("L77").Value = 61391.39999999999
$ws.Range("N77").Value = -70751.39999999999
$ws.Range("H135").Value = 74555
$ws.Range("J135").Value = 74555
$ws.Range("L135").Value = 74555
$ws.Range("N135").Value = -84695
$ws.Range("H136").Value = 3088.1777
$ws.Range("I136").Value = 2354.8948
$ws.Range("K136").Value = 7064.6844
$ws.Range("M136").Value = -4514.6844
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = ""
